$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.96"
$ws.Range("E2").Value = "'-0.09%"
$ws.Range("D3").Value = "'28.29"
$ws.Range("D4").Value = "'5.293"
$ws.Range("E4").Value = "'2.15%"
$ws.Range("D5").Value = "'0.05709"
$ws.Range("E5").Value = "'-0.48%"
$ws.Range("D6").Value = "'6.643"
$ws.Range("E6").Value = "'1.15%"
$ws.Range("D7").Value = "'3.216"
$ws.Range("E7").Value = "'3.48%"
$ws.Range("D8").Value = "'0.8628"
$ws.Range("E8").Value = "'0.44%"
$ws.Range("D9").Value = "'0.8845"
$ws.Range("E9").Value = "'2.44%"
$ws.Range("D10").Value = "'0.1389"
$ws.Range("E10").Value = "'1.90%"
$ws.Range("D11").Value = "'0.07089"
$ws.Range("E11").Value = "'0.03%"
$ws.Range("D12").Value = "'0.03152"
$ws.Range("E12").Value = "'3.07%"
$ws.Range("D13").Value = "'0.09234"
$ws.Range("E13").Value = "'-1.46%"
$ws.Range("D14").Value = "'0.001527"
$ws.Range("E14").Value = "'-0.75%"
$ws.Range("D15").Value = "'0.0005959"
$ws.Range("E15").Value = "'-94.20%"
$ws.Range("D16").Value = "'0.006023"
$ws.Range("E16").Value = "'0.51%"
$ws.Range("D17").Value = "'3.496"
$ws.Range("E17").Value = "'0.18%"
$ws.Range("E18").Value = "'-4.58%"
$ws.Range("D19").Value = "'0.3166"
$ws.Range("E19").Value = "'-0.92%"
$ws.Range("D20").Value = "'0.03350"
$ws.Range("E20").Value = "'1.12%"
$ws.Range("D22").Value = "'3.488"
$ws.Range("E22").Value = "'-0.16%"
$ws.Range("D23").Value = "'0.04101"
$ws.Range("E23").Value = "'-0.95%"
$ws.Range("D25").Value = "'0.001219"
$ws.Range("E25").Value = "'-0.69%"
$ws.Range("D26").Value = "'0.004164"
$ws.Range("E26").Value = "'-16.59%"
$ws.Range("E27").Value = "'-0.92%"
$ws.Range("D28").Value = "'0.0001445"
$ws.Range("E40").Value = "'1.06%"
$ws.Range("D41").Value = "'0.1068"
$ws.Range("E41").Value = "'-0.21%"
$ws.Range("D42").Value = "'0.002199"
$ws.Range("E42").Value = "'4.67%"
$ws.Range("E43").Value = "'-49.54%"
$ws.Range("D44").Value = "'0.009476"
$ws.Range("E44").Value = "'11.92%"
$ws.Range("D45").Value = "'0.00005268"
$ws.Range("E45").Value = "'-0.29%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.09%"
$ws.Range("D47").Value = "'0.08908"
$ws.Range("E47").Value = "'56.25%"
$ws.Range("E48").Value = "'-0.21%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.09%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.09%"
